$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dane")

$ws.Range("A2").Value = "Dolnoslaskie"
$ws.Range("A6").Value = "Lodzkie"
$ws.Range("A7").Value = "Malopolskie"
$ws.Range("A13").Value = "Slaskie"
$ws.Range("A14").Value = "Swietokrzyskie"
$ws.Range("A15").Value = "Warminsko-mazurskie"

$ws.Range("A16").Select()
